$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in missing grades on row 8 (student #4): D8:G8 were empty cells,
# H8:K8 existed but were blank. All become 5, keeping existing cell styles.
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 5
$ws.Range("G8").Value = 5
$ws.Range("H8").Value = 5
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 5
$ws.Range("K8").Value = 5

# Fill in missing grade on row 14 (student #10)
$ws.Range("H14").Value = 5

# Move the view back to the top (unfreeze scroll) and select L8
$ws.Range("L8").Select()
